# Add a new handed-back file ("d9261fb8-c2c4-4bcc-a3a4-3d2221a0094a") as
# row 9 to the Overview / zh-cn / de-de sheets, mirroring the existing
# "Ready for handoff" rows (e.g. c647c878-740c-4502-98a1-6a9c3136ac42).
#
# Notes:
#  - A leading "'" forces Excel to store a literal empty string (instead
#    of deleting the cell) or literal text "True"/"False" (instead of
#    auto-coercing to a boolean cell).

$wb = $excel.ActiveWorkbook

$guid = "d9261fb8-c2c4-4bcc-a3a4-3d2221a0094a"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = "$guid.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0/e2e/$guid.md",
    "",
    "",
    "e2e\$guid.md"
)
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = "'"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-31 20:56:14"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "'False"
$wsZhCn.Range("G9").Value = "$guid.0cca7d5e71214a1ff9c206d4b6f2a020cf8f3325.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-31 20:55:58"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = "'"
$wsZhCn.Range("J9").Value = "'"
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = "'"
$wsZhCn.Range("M9").Value = "'True"
$wsZhCn.Range("N9").Value = "'"
$wsZhCn.Range("O9").Value = "'False"
$wsZhCn.Range("P9").Value = "'"

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P9"))

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0c1d2/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "'False"
$wsDeDe.Range("G9").Value = "$guid.0cca7d5e71214a1ff9c206d4b6f2a020cf8f3325.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-31 20:56:14"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = "'"
$wsDeDe.Range("J9").Value = "'"
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = "'"
$wsDeDe.Range("M9").Value = "'True"
$wsDeDe.Range("N9").Value = "'"
$wsDeDe.Range("O9").Value = "'False"
$wsDeDe.Range("P9").Value = "'"

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P9"))
